$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string but must remain TEXT
# (matches the source data which stores prices/etc. as inline strings).
# Marking the cell as Text before assignment prevents Excel from silently
# converting these into numbers (which would also lose formatting like
# trailing zeros, e.g. "246.00" -> 246).
$textForcedCells = @(
    'D2'
    'D3'
    'D4'
    'D5'
    'D6'
    'D7'
    'D8'
    'D9'
    'D10'
    'D11'
    'D12'
    'D13'
    'D14'
    'D15'
    'D16'
    'D17'
    'D18'
    'D19'
    'D20'
    'D21'
    'D24'
    'D25'
    'D26'
    'D27'
    'D40'
    'D41'
    'D42'
    'D43'
    'D44'
    'D45'
    'D47'
    'D48'
    'D49'
)
foreach ($ref in $textForcedCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the new cell values.
$ws.Range('D2').Value = '246.00'
$ws.Range('D3').Value = '21.98'
$ws.Range('D4').Value = '5.366'
$ws.Range('D5').Value = '0.05854'
$ws.Range('D6').Value = '3.395'
$ws.Range('D7').Value = '6.375'
$ws.Range('D8').Value = '0.8147'
$ws.Range('D9').Value = '1.021'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '0.1418'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '0.04175'
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.07357'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.02985'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'MCDex'
$ws.Range('C14').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D14').Value = '4.220'
$ws.Range('E14').Value = '13MCDexMCB'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = '0.09394'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = '0.001590'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = '0.04802'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').Value = '0.0005891'
$ws.Range('E18').Value = '17OneONE'
$ws.Range('D19').Value = '0.006008'
$ws.Range('D20').Value = '0.004083'
$ws.Range('D21').Value = '0.0009818'
$ws.Range('D24').Value = '2.232'
$ws.Range('D25').Value = '0.3241'
$ws.Range('D26').Value = '0.1295'
$ws.Range('D27').Value = '0.0002484'
$ws.Range('D40').Value = '0.03858'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = '0.1073'
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = '0.002410'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').Value = '0.003039'
$ws.Range('E43').Value = '42KickTokenKICKWorstin24h'
$ws.Range('D44').Value = '0.005070'
$ws.Range('D45').Value = '0.00005629'
$ws.Range('D47').Value = '0.8002'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOINBestin24h'
$ws.Range('D48').Value = '0.09273'
$ws.Range('E48').Value = '47BOLOBOLO'
$ws.Range('D49').Value = '0.00002101'
